$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Primary Key" column header next to "Table for Books" title
$ws.Range("B1").Value = "Primary Key"

# Update Price's data type from "double" to "float8"
$ws.Range("F3").Value = "float8"

# Resize columns to fit the new content
$ws.Columns.Item(2).ColumnWidth = 10.666666666666668
$ws.Columns.Item(8).ColumnWidth = 15.833333333333332

# Update the active selection/cursor position
[void]$ws.Range("E15").Select()
